# Auto-generated edit script applying cell value updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Addr='D2'; Val='26.536.43'; Numeric=$false},
    @{Addr='E2'; Val='  +4.03%  '; Numeric=$false},
    @{Addr='D3'; Val='1.738.62'; Numeric=$false},
    @{Addr='E3'; Val='  +4.40%  '; Numeric=$false},
    @{Addr='D4'; Val='0.9995'; Numeric=$true},
    @{Addr='E4'; Val='  +0.11%  '; Numeric=$false},
    @{Addr='D5'; Val='245.66'; Numeric=$true},
    @{Addr='E5'; Val='  +4.68%  '; Numeric=$false},
    @{Addr='D6'; Val='1.000'; Numeric=$true},
    @{Addr='E6'; Val='  +0.03%  '; Numeric=$false},
    @{Addr='D7'; Val='0.4799'; Numeric=$true},
    @{Addr='E7'; Val='  +3.12%  '; Numeric=$false},
    @{Addr='D8'; Val='0.2683'; Numeric=$true},
    @{Addr='E8'; Val='  +4.08%  '; Numeric=$false},
    @{Addr='D9'; Val='0.06244'; Numeric=$true},
    @{Addr='E9'; Val='  +1.68%  '; Numeric=$false},
    @{Addr='D10'; Val='1.738.67'; Numeric=$false},
    @{Addr='E10'; Val='  +4.50%  '; Numeric=$false},
    @{Addr='D11'; Val='0.07125'; Numeric=$true},
    @{Addr='E11'; Val='  +2.57%  '; Numeric=$false},
    @{Addr='D12'; Val='15.81'; Numeric=$true},
    @{Addr='E12'; Val='  +7.86%  '; Numeric=$false},
    @{Addr='D13'; Val='0.6203'; Numeric=$true},
    @{Addr='E13'; Val='  +8.50%  '; Numeric=$false},
    @{Addr='D14'; Val='4.542'; Numeric=$true},
    @{Addr='E14'; Val='  +3.88%  '; Numeric=$false},
    @{Addr='D15'; Val='77.15'; Numeric=$true},
    @{Addr='E15'; Val='  +2.81%  '; Numeric=$false},
    @{Addr='D16'; Val='0.9999'; Numeric=$true},
    @{Addr='E16'; Val='  +0.00%  '; Numeric=$false},
    @{Addr='D17'; Val='26.541.53'; Numeric=$false},
    @{Addr='E17'; Val='  +4.08%  '; Numeric=$false},
    @{Addr='D18'; Val='1.000'; Numeric=$true},
    @{Addr='B19'; Val='Avalanche'; Numeric=$false},
    @{Addr='C19'; Val='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Numeric=$false},
    @{Addr='D19'; Val='11.78'; Numeric=$true},
    @{Addr='E19'; Val='  +3.56%  '; Numeric=$false},
    @{Addr='B20'; Val='ShibaInu'; Numeric=$false},
    @{Addr='C20'; Val='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Numeric=$false},
    @{Addr='D20'; Val='0.000006893'; Numeric=$true},
    @{Addr='E20'; Val='  +2.87%  '; Numeric=$false},
    @{Addr='D21'; Val='1.961.93'; Numeric=$false},
    @{Addr='E21'; Val='  +4.46%  '; Numeric=$false},
    @{Addr='D22'; Val='4.607'; Numeric=$true},
    @{Addr='E22'; Val='  +3.88%  '; Numeric=$false},
    @{Addr='D23'; Val='8.896'; Numeric=$true},
    @{Addr='E23'; Val='  +2.01%  '; Numeric=$false},
    @{Addr='D24'; Val='5.350'; Numeric=$true},
    @{Addr='E24'; Val='  +2.51%  '; Numeric=$false},
    @{Addr='D25'; Val='135.90'; Numeric=$true},
    @{Addr='E25'; Val='  +0.09%  '; Numeric=$false},
    @{Addr='D26'; Val='15.37'; Numeric=$true},
    @{Addr='E26'; Val='  +3.53%  '; Numeric=$false},
    @{Addr='D27'; Val='1.808'; Numeric=$true},
    @{Addr='E27'; Val='  +6.23%  '; Numeric=$false},
    @{Addr='D28'; Val='1.421'; Numeric=$true},
    @{Addr='E28'; Val='  +4.29%  '; Numeric=$false},
    @{Addr='E29'; Val='  +3.04%  '; Numeric=$false},
    @{Addr='D30'; Val='4.002'; Numeric=$true},
    @{Addr='E30'; Val='  +1.50%  '; Numeric=$false},
    @{Addr='D31'; Val='3.736'; Numeric=$true},
    @{Addr='E31'; Val='  +3.58%  '; Numeric=$false},
    @{Addr='D32'; Val='0.07864'; Numeric=$true},
    @{Addr='E32'; Val='  +1.98%  '; Numeric=$false},
    @{Addr='D33'; Val='0.04583'; Numeric=$true},
    @{Addr='E33'; Val='  +6.57%  '; Numeric=$false},
    @{Addr='B34'; Val='Frax'; Numeric=$false},
    @{Addr='C34'; Val='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; Numeric=$false},
    @{Addr='D34'; Val='0.9993'; Numeric=$true},
    @{Addr='E34'; Val='  +0.05%  '; Numeric=$false},
    @{Addr='B35'; Val='HuobiToken'; Numeric=$false},
    @{Addr='C35'; Val='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Numeric=$false},
    @{Addr='D35'; Val='2.615'; Numeric=$true},
    @{Addr='E35'; Val='  -0.18%  '; Numeric=$false},
    @{Addr='B36'; Val='ImmutableX'; Numeric=$false},
    @{Addr='C36'; Val='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Numeric=$false},
    @{Addr='D36'; Val='0.6385'; Numeric=$true},
    @{Addr='E36'; Val='  +6.40%  '; Numeric=$false},
    @{Addr='B37'; Val='ARBITRUM'; Numeric=$false},
    @{Addr='C37'; Val='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Numeric=$false},
    @{Addr='D37'; Val='1.001'; Numeric=$true},
    @{Addr='E37'; Val='  +6.04%  '; Numeric=$false},
    @{Addr='B38'; Val='TrustWalletToken'; Numeric=$false},
    @{Addr='C38'; Val='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Numeric=$false},
    @{Addr='D38'; Val='0.9308'; Numeric=$true},
    @{Addr='E38'; Val='  +0.79%  '; Numeric=$false},
    @{Addr='B39'; Val='Quant'; Numeric=$false},
    @{Addr='C39'; Val='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Numeric=$false},
    @{Addr='D39'; Val='112.84'; Numeric=$true},
    @{Addr='E39'; Val='  +8.80%  '; Numeric=$false},
    @{Addr='B40'; Val='MXToken'; Numeric=$false},
    @{Addr='C40'; Val='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Numeric=$false},
    @{Addr='D40'; Val='2.427'; Numeric=$true},
    @{Addr='E40'; Val='  -2.05%  '; Numeric=$false},
    @{Addr='B41'; Val='RenderToken'; Numeric=$false},
    @{Addr='C41'; Val='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Numeric=$false},
    @{Addr='D41'; Val='1.979'; Numeric=$true},
    @{Addr='E41'; Val='  +8.88%  '; Numeric=$false},
    @{Addr='B42'; Val='PaxDollar'; Numeric=$false},
    @{Addr='C42'; Val='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; Numeric=$false},
    @{Addr='D42'; Val='1.003'; Numeric=$true},
    @{Addr='E42'; Val='  +0.35%  '; Numeric=$false},
    @{Addr='B43'; Val='FraxShare'; Numeric=$false},
    @{Addr='C43'; Val='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Numeric=$false},
    @{Addr='D43'; Val='5.743'; Numeric=$true},
    @{Addr='E43'; Val='  +14.18%  '; Numeric=$false},
    @{Addr='B44'; Val='VeChain'; Numeric=$false},
    @{Addr='C44'; Val='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Numeric=$false},
    @{Addr='D44'; Val='0.01511'; Numeric=$true},
    @{Addr='E44'; Val='  +3.37%  '; Numeric=$false},
    @{Addr='B45'; Val='TheSandbox'; Numeric=$false},
    @{Addr='C45'; Val='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Numeric=$false},
    @{Addr='D45'; Val='0.3908'; Numeric=$true},
    @{Addr='E45'; Val='  +5.39%  '; Numeric=$false},
    @{Addr='B46'; Val='Aptos'; Numeric=$false},
    @{Addr='C46'; Val='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Numeric=$false},
    @{Addr='D46'; Val='6.900'; Numeric=$true},
    @{Addr='E46'; Val='  +12.69%  '; Numeric=$false},
    @{Addr='B47'; Val='Algorand'; Numeric=$false},
    @{Addr='C47'; Val='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Numeric=$false},
    @{Addr='D47'; Val='0.1206'; Numeric=$true},
    @{Addr='E47'; Val='  +9.08%  '; Numeric=$false},
    @{Addr='B48'; Val='Cronos'; Numeric=$false},
    @{Addr='C48'; Val='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Numeric=$false},
    @{Addr='D48'; Val='0.05335'; Numeric=$true},
    @{Addr='E48'; Val='  +1.33%  '; Numeric=$false},
    @{Addr='D49'; Val='7.917'; Numeric=$true},
    @{Addr='E49'; Val='  +6.79%  '; Numeric=$false},
    @{Addr='B50'; Val='Elrond'; Numeric=$false},
    @{Addr='C50'; Val='https://coinranking.com/coin/omwkOTglq+elrond-egld'; Numeric=$false},
    @{Addr='D50'; Val='30.73'; Numeric=$true},
    @{Addr='E50'; Val='  +3.43%  '; Numeric=$false},
    @{Addr='B51'; Val='NEARProtocol'; Numeric=$false},
    @{Addr='C51'; Val='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Numeric=$false},
    @{Addr='D51'; Val='1.256'; Numeric=$true},
    @{Addr='E51'; Val='  +5.31%  '; Numeric=$false},
)

foreach ($item in $changes) {
    $r = $ws.Range($item.Addr)
    if ($item.Numeric) {
        # Force text storage so numeric-looking strings keep their exact
        # textual representation (e.g. trailing zeros, dot-thousand separators)
        $r.NumberFormat = '@'
        $r.Value = $item.Val
        $r.Style = 'Normal'
    } else {
        $r.Value = $item.Val
    }
}

Write-Host ('Applied ' + $changes.Count + ' cell updates')
